{"js": "// Replace each two-digit multiplication prompt in the table with its new\n// value. Every prompt text in this document is unique, so searching for the\n// exact old text and replacing it with the new text (in document order)\n// reproduces the diff unambiguously.\nconst replacements = [\n  [\"94\u00d712=\", \"38\u00d788=\"],\n  [\"88\u00d717=\", \"68\u00d725=\"],\n  [\"20\u00d738=\", \"59\u00d754=\"],\n  [\"46\u00d760=\", \"77\u00d761=\"],\n  [\"19\u00d795=\", \"83\u00d719=\"],\n  [\"53\u00d761=\", \"20\u00d795=\"],\n  [\"41\u00d794=\", \"34\u00d760=\"],\n  [\"23\u00d783=\", \"63\u00d778=\"],\n  [\"49\u00d724=\", \"87\u00d736=\"],\n  [\"54\u00d718=\", \"71\u00d740=\"],\n  [\"40\u00d736=\", \"27\u00d749=\"],\n  [\"17\u00d797=\", \"16\u00d716=\"],\n  [\"74\u00d720=\", \"66\u00d784=\"],\n  [\"74\u00d747=\", \"59\u00d712=\"],\n  [\"56\u00d761=\", \"71\u00d793=\"],\n  [\"82\u00d791=\", \"27\u00d751=\"],\n  [\"77\u00d776=\", \"59\u00d790=\"],\n  [\"92\u00d795=\", \"54\u00d740=\"],\n  [\"25\u00d752=\", \"49\u00d772=\"],\n  [\"68\u00d790=\", \"94\u00d733=\"],\n  [\"67\u00d768=\", \"89\u00d759=\"],\n  [\"84\u00d726=\", \"25\u00d735=\"],\n  [\"84\u00d739=\", \"69\u00d770=\"],\n  [\"60\u00d792=\", \"42\u00d782=\"],\n  [\"58\u00d731=\", \"85\u00d794=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication prompt in the table with its new\n# value. Every prompt text in this document is unique, so a Find/Replace\n# pass for each exact old -> new pair (run once, wdReplaceAll) reproduces\n# the diff unambiguously.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"94\u00d712=\", \"38\u00d788=\"),\n  @(\"88\u00d717=\", \"68\u00d725=\"),\n  @(\"20\u00d738=\", \"59\u00d754=\"),\n  @(\"46\u00d760=\", \"77\u00d761=\"),\n  @(\"19\u00d795=\", \"83\u00d719=\"),\n  @(\"53\u00d761=\", \"20\u00d795=\"),\n  @(\"41\u00d794=\", \"34\u00d760=\"),\n  @(\"23\u00d783=\", \"63\u00d778=\"),\n  @(\"49\u00d724=\", \"87\u00d736=\"),\n  @(\"54\u00d718=\", \"71\u00d740=\"),\n  @(\"40\u00d736=\", \"27\u00d749=\"),\n  @(\"17\u00d797=\", \"16\u00d716=\"),\n  @(\"74\u00d720=\", \"66\u00d784=\"),\n  @(\"74\u00d747=\", \"59\u00d712=\"),\n  @(\"56\u00d761=\", \"71\u00d793=\"),\n  @(\"82\u00d791=\", \"27\u00d751=\"),\n  @(\"77\u00d776=\", \"59\u00d790=\"),\n  @(\"92\u00d795=\", \"54\u00d740=\"),\n  @(\"25\u00d752=\", \"49\u00d772=\"),\n  @(\"68\u00d790=\", \"94\u00d733=\"),\n  @(\"67\u00d768=\", \"89\u00d759=\"),\n  @(\"84\u00d726=\", \"25\u00d735=\"),\n  @(\"84\u00d739=\", \"69\u00d770=\"),\n  @(\"60\u00d792=\", \"42\u00d782=\"),\n  @(\"58\u00d731=\", \"85\u00d794=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
